$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 12 - this shifts every row below it (including
# the data rows 12-15 and the formatting-only "spacer" rows further down the
# sheet, e.g. 47->48, 54->55, ... 175->176) down by one.
$null = $ws.Rows("12:12").Insert(1)

# The new row is blank after Insert; pick up the same visual formatting
# (borders/alignment/number-format) used by the rest of the table by copying
# it from the row right below (which used to be row 12 before the shift).
$ws.Range("A13:E13").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new test case's data.
$ws.Range("A12").Value = "MultipleErrorTypes"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "Finished"
$ws.Range("E12").Value = "Contains 1 partially automated test case."

# The single-cell conditional formatting rule-set that used to watch D12 now
# needs to watch D13 (the row that used to be 12, now pushed one row down) -
# both the range it applies to and the formula text that references the cell.
# (Find every matching rule first: re-pointing one of them re-points the
# whole sqref group, which would make the later items stop matching.)
$fcs = $ws.Cells.FormatConditions
$rulesToFix = @()
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.AppliesTo.Address() -eq '$D$12') {
        $rulesToFix += $fc
    }
}
foreach ($fc in $rulesToFix) {
    $newFormula = $fc.Formula1 -replace 'D12\b', 'D13'
    $fc.ModifyAppliesToRange($ws.Range("D13"))
    $fc.Formula1 = $newFormula
}

# Leave the selection where the author's session ended up.
$null = $ws.Range("E13").Select()
